# Apply odds updates to the "Jogos da Semana" sheet for rows 10-13
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "H10" = 3.1
    "I10" = 2.3
    "J10" = 3.55
    "L10" = 2.8
    "N10" = 2.77
    "R10" = 1.25
    "W10" = 8.5
    "X10" = 15
    "AA10" = 28
    "AB10" = 37
    "AC10" = 8.5
    "AH10" = 7.6
    "AI10" = 11.25
    "AK10" = 24
    "AL10" = 19

    "O11" = 1.62
    "P11" = 2.2
    "Q11" = 2.5
    "R11" = 1.5

    "I12" = 2.88
    "K12" = 2.25
    "O12" = 1.58
    "S12" = 1.33
    "T12" = 3.25
    "U12" = 1.57
    "V12" = 2.25
    "AB12" = 23
    "AJ12" = 11
    "AN12" = 1.03
    "AO12" = 10

    "G13" = 2.5
    "I13" = 2.63
    "J13" = 3.1
    "L13" = 3.25
    "O13" = 2.1
    "P13" = 1.7
    "Q13" = 3.75
    "R13" = 1.25
    "U13" = 1.8
    "V13" = 1.91
    "W13" = 8.5
    "X13" = 12
    "Z13" = 23
    "AA13" = 21
    "AB13" = 34
    "AH13" = 9
    "AI13" = 13
    "AJ13" = 11
    "AK13" = 26
    "AL13" = 23
    "AN13" = 1.07
    "AO13" = 7.5
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
